$d = $word.ActiveDocument

$d.Content.Find.Execute("54÷6=9, 0", $true, $false, $false, $false, $false, $true, 1, $false, "57÷7=8, 1", 2) | Out-Null
$d.Content.Find.Execute("99÷4=24, 3", $true, $false, $false, $false, $false, $true, 1, $false, "72÷6=12, 0", 2) | Out-Null
$d.Content.Find.Execute("83÷5=16, 3", $true, $false, $false, $false, $false, $true, 1, $false, "18÷3=6, 0", 2) | Out-Null
$d.Content.Find.Execute("71÷4=17, 3", $true, $false, $false, $false, $false, $true, 1, $false, "60÷6=10, 0", 2) | Out-Null
$d.Content.Find.Execute("77÷4=19, 1", $true, $false, $false, $false, $false, $true, 1, $false, "65÷6=10, 5", 2) | Out-Null
$d.Content.Find.Execute("10÷4=2, 2", $true, $false, $false, $false, $false, $true, 1, $false, "57÷6=9, 3", 2) | Out-Null
$d.Content.Find.Execute("27÷2=13, 1", $true, $false, $false, $false, $false, $true, 1, $false, "14÷2=7, 0", 2) | Out-Null
$d.Content.Find.Execute("76÷7=10, 6", $true, $false, $false, $false, $false, $true, 1, $false, "52÷9=5, 7", 2) | Out-Null
$d.Content.Find.Execute("15÷8=1, 7", $true, $false, $false, $false, $false, $true, 1, $false, "75÷7=10, 5", 2) | Out-Null
$d.Content.Find.Execute("64÷5=12, 4", $true, $false, $false, $false, $false, $true, 1, $false, "94÷3=31, 1", 2) | Out-Null
$d.Content.Find.Execute("89÷7=12, 5", $true, $false, $false, $false, $false, $true, 1, $false, "16÷6=2, 4", 2) | Out-Null
$d.Content.Find.Execute("50÷2=25, 0", $true, $false, $false, $false, $false, $true, 1, $false, "80÷3=26, 2", 2) | Out-Null
$d.Content.Find.Execute("29÷5=5, 4", $true, $false, $false, $false, $false, $true, 1, $false, "65÷6=10, 5", 2) | Out-Null
$d.Content.Find.Execute("75÷9=8, 3", $true, $false, $false, $false, $false, $true, 1, $false, "39÷8=4, 7", 2) | Out-Null
$d.Content.Find.Execute("28÷3=9, 1", $true, $false, $false, $false, $false, $true, 1, $false, "67÷5=13, 2", 2) | Out-Null
$d.Content.Find.Execute("21÷4=5, 1", $true, $false, $false, $false, $false, $true, 1, $false, "41÷7=5, 6", 2) | Out-Null
$d.Content.Find.Execute("99÷8=12, 3", $true, $false, $false, $false, $false, $true, 1, $false, "28÷2=14, 0", 2) | Out-Null
$d.Content.Find.Execute("29÷7=4, 1", $true, $false, $false, $false, $false, $true, 1, $false, "53÷3=17, 2", 2) | Out-Null
$d.Content.Find.Execute("68÷7=9, 5", $true, $false, $false, $false, $false, $true, 1, $false, "90÷5=18, 0", 2) | Out-Null
$d.Content.Find.Execute("58÷4=14, 2", $true, $false, $false, $false, $false, $true, 1, $false, "88÷5=17, 3", 2) | Out-Null
$d.Content.Find.Execute("59÷8=7, 3", $true, $false, $false, $false, $false, $true, 1, $false, "17÷9=1, 8", 2) | Out-Null
$d.Content.Find.Execute("34÷9=3, 7", $true, $false, $false, $false, $false, $true, 1, $false, "30÷3=10, 0", 2) | Out-Null
$d.Content.Find.Execute("69÷4=17, 1", $true, $false, $false, $false, $false, $true, 1, $false, "97÷9=10, 7", 2) | Out-Null
$d.Content.Find.Execute("51÷6=8, 3", $true, $false, $false, $false, $false, $true, 1, $false, "10÷4=2, 2", 2) | Out-Null
$d.Content.Find.Execute("44÷5=8, 4", $true, $false, $false, $false, $false, $true, 1, $false, "34÷9=3, 7", 2) | Out-Null
